$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.472.35'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.604.42'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.92'
$ws.Range('E5').Value = '  +1.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.65'
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('E8').Value = '  +0.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.630.05'
$ws.Range('E9').Value = '  +1.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.50'
$ws.Range('E10').Value = '  -2.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.105'
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('E12').Value = '  -4.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.364'
$ws.Range('E13').Value = '  +1.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.067.31'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.481.38'
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.21'
$ws.Range('E16').Value = '  -0.38%  '
$ws.Range('E17').Value = '  +2.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.614.17'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.35'
$ws.Range('E19').Value = '  +8.83%  '
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '346.98'
$ws.Range('E21').Value = '  +2.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.96'
$ws.Range('E22').Value = '  +6.68%  '
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.525'
$ws.Range('E24').Value = '  +10.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.22'
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  -1.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.73'
$ws.Range('E28').Value = '  +3.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0788'
$ws.Range('E29').Value = '  +0.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.84'
$ws.Range('E30').Value = '  +9.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.37'
$ws.Range('E31').Value = '  +2.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '161.19'
$ws.Range('E33').Value = '  +1.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.50'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.22'
$ws.Range('E35').Value = '  +3.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.979'
$ws.Range('E36').Value = '  +8.53%  '
$ws.Range('E37').Value = '  +3.66%  '
$ws.Range('E38').Value = '  +7.33%  '
$ws.Range('E39').Value = '  +1.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.84'
$ws.Range('E40').Value = '  +4.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.843'
$ws.Range('E41').Value = '  -3.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '294.54'
$ws.Range('E42').Value = '  -0.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '137.25'
$ws.Range('E43').Value = '  -2.13%  '
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('E45').Value = '  +2.12%  '
$ws.Range('E46').Value = '  +0.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.77'
$ws.Range('E47').Value = '  +2.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0544'
$ws.Range('E48').Value = '  +2.12%  '
$ws.Range('E49').Value = '  +8.46%  '
$ws.Range('E50').Value = '  +1.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.72'
$ws.Range('E51').Value = '  +0.81%  '
